# Fix the "weird bug": the sheet was missing a "Date of Last Update" column,
# which caused every column from "Shipping-Related Issues" onward to carry
# the wrong header/value (everything was shifted one slot to the left). Also
# a duplicate "TicketID" column (carried over from the source export) needs
# to be restored right after "Time Worked (Minutes)".
#
# Net effect: insert two new columns at I:J.
#   I = "Date of Last Update"  (the value that used to incorrectly live
#        under "Shipping-Related Issues")
#   J = duplicate "TicketID" column
# and correct the now-shifted "Shipping-Related Issues" cell, which had
# been showing the date instead of the real issue text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns in front of the old column I
# ("Shipping-Related Issues" / "2017-08-15T09:31:24.747000").
$ws.Columns("I:J").Insert()

# New column I: "Date of Last Update" -- reuse the date value that used to
# (incorrectly) sit in the old "Shipping-Related Issues" column, now at K.
$ws.Range("I1").Value2 = "Date of Last Update"
$ws.Range("K2").Copy()
$ws.Range("I2").PasteSpecial(-4163)

# New column J: duplicate of the "TicketID" column (A).
$ws.Range("A1").Copy()
$ws.Range("J1").PasteSpecial(-4163)
$ws.Range("A2").Copy()
$ws.Range("J2").PasteSpecial(-4163)

# The (now shifted) "Shipping-Related Issues" cell still holds the old date
# value -- replace it with the real issue text.
$ws.Range("K2").Value2 = "Incorrect Parts Shipped"
